$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 gets a new text value
$ws.Range("A1").Value = "Nuevo valor"

# New row of data on row 3 (B3:D3)
$ws.Range("B3").Value = "nuevaFila"
$ws.Range("C3").Value = "valor1"
$ws.Range("D3").Value = "valor2"
